$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells whose new text value looks like a plain number need to be
# forced back to text (NumberFormat "@") so Excel does not silently convert
# them into a numeric cell; Style is then reset to "Normal" so no stray
# cell-level style (s="...") is left behind on the cell.

$ws.Range("D2").Value = "57.735.64"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "2.440.70"
$ws.Range("E3").Value = "  -2.93%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "

$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  +0.86%  "

$ws.Range("D9").Value = "2.440.59"
$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.60%  "

$ws.Range("E11").Value = "  -2.31%  "

$ws.Range("E12").Value = "  -4.26%  "

$ws.Range("E13").Value = "  -2.38%  "

$ws.Range("D14").Value = "2.874.25"
$ws.Range("E14").Value = "  -2.29%  "

$ws.Range("D15").Value = "57.670.64"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.34%  "

$ws.Range("D18").Value = "2.439.60"
$ws.Range("E18").Value = "  -2.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.15%  "

$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.44%  "

$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.403"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "

$ws.Range("E26").Value = "  +0.51%  "

$ws.Range("E27").Value = "  -1.50%  "

$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.79%  "

$ws.Range("D30").Value = "0.0₃0733"
$ws.Range("E30").Value = "  -2.62%  "

$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("E32").Value = "  -2.51%  "

$ws.Range("E33").Value = "  -4.20%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "

$ws.Range("E37").Value = "  -4.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  -1.69%  "

$ws.Range("E41").Value = "  +2.75%  "

$ws.Range("E42").Value = "  -1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "262.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.49%  "

$ws.Range("E44").Value = "  -2.33%  "

$ws.Range("E45").Value = "  -4.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0922"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0492"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.56%  "

$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.64%  "
